# "correct text wrap length" - trims the extra test rows that were produced
# by a runaway test script, fixes the mis-typed header + the first
# date/time sample, and re-tunes the column widths now that the sheet is
# back to a sane size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel's
# autocomplete/date-sniffing kicking in (a plain .Value = "06-11-2023"
# gets silently reinterpreted as a date serial). Route it through a
# text-formula in a scratch cell, then Paste-Special "values only" so the
# destination ends up with a true shared-string / General-format cell,
# exactly like the original file's cells.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# --- drop the extra test rows (old rows 10-18) ------------------------
$ws.Rows("10:18").Delete()

# --- header fix ---------------------------------------------------------
Set-TextValue $ws.Range("D2") "id"

# --- fix the sample date (rows 3-9 all shared the same placeholder date) -
Set-TextValue $ws.Range("A3") "06-11-2023"
Set-TextValue $ws.Range("A4") "06-11-2023"
Set-TextValue $ws.Range("A5") "06-11-2023"
Set-TextValue $ws.Range("A6") "06-11-2023"
Set-TextValue $ws.Range("A7") "06-11-2023"
Set-TextValue $ws.Range("A8") "06-11-2023"
Set-TextValue $ws.Range("A9") "06-11-2023"

# --- fix the sample times, row by row -----------------------------------
Set-TextValue $ws.Range("B3") "15:36:11"
Set-TextValue $ws.Range("B4") "15:36:45"
Set-TextValue $ws.Range("B5") "15:36:46"
Set-TextValue $ws.Range("B6") "15:36:47"
Set-TextValue $ws.Range("B7") "15:36:48"
Set-TextValue $ws.Range("B8") "15:36:49"
Set-TextValue $ws.Range("B9") "15:36:50"

# --- re-apply the (now explicit) wrapped-row height ---------------------
$ws.Rows.Item(3).RowHeight = 28.8

# --- column widths: COM's ColumnWidth stores 5/6 of a character wider
# than the number you set (Excel's default-font cell-padding), so back
# that off to land on the exact widths used by the corrected layout.
$pad = 5/6
$ws.Columns.Item(1).ColumnWidth = 15 - $pad
$ws.Columns("B:D").ColumnWidth = 13 - $pad
$ws.Columns.Item(5).ColumnWidth = 40 - $pad
$ws.Columns.Item(6).ColumnWidth = 15 - $pad
$ws.Columns.Item(7).ColumnWidth = 13 - $pad
$ws.Columns.Item(8).ColumnWidth = 40 - $pad

# --- restore the view: frozen header rows, selection on H3 --------------
$ws.Range("H3").Select()
